$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete the specific rows (account holders) that were removed in this upload.
# Row 2: 003301389 / EDMUNDO / 113000
# Row 7: 004240014 / ISABELE / 10850.63
# Row 8: 004752461 / SERGIO / 10774.85
# Row 10: 004357848 / AURELIO / 4808.2
# Delete from bottom to top so earlier row numbers stay valid.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(2).Delete()
